$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 23288.47245586674
$ws.Cells.Item(2, 3).Value = 38.15225880359203
$ws.Cells.Item(2, 4).Value = 0.6628867188133686

$ws.Cells.Item(3, 2).Value = 23884.15661053462
$ws.Cells.Item(3, 3).Value = 49.29702604143551
$ws.Cells.Item(3, 4).Value = 0.6733790583472918

$ws.Cells.Item(4, 2).Value = 24806.89636359204
$ws.Cells.Item(4, 3).Value = 219.6108457979281
$ws.Cells.Item(4, 4).Value = 0.6888162620334709

$ws.Cells.Item(5, 2).Value = 25122.2562531207
$ws.Cells.Item(5, 3).Value = 83.75949344712977
$ws.Cells.Item(5, 4).Value = 0.6938726693684723

$ws.Cells.Item(6, 1).Value = "DM0 Subgroup 1"
$ws.Cells.Item(6, 2).Value = 52896.85629800297
$ws.Cells.Item(6, 3).Value = 42944.76085826926
$ws.Cells.Item(6, 4).Value = 0.7213184675172366

$ws.Cells.Item(7, 1).Value = "DM1 Subgroup 1"
$ws.Cells.Item(7, 2).Value = 53288.81922628694
$ws.Cells.Item(7, 3).Value = 25184.58254266901
$ws.Cells.Item(7, 4).Value = 0.7240632928562425

$ws.Cells.Item(8, 1).Value = "DM2 Subgroup 1"
$ws.Cells.Item(8, 2).Value = 53421.74929117745
$ws.Cells.Item(8, 3).Value = 5907.321406219426
$ws.Cells.Item(8, 4).Value = 0.7249887953683581

$ws.Cells.Item(9, 2).Value = 13816.13394493324
$ws.Cells.Item(9, 3).Value = 3.336464188559541
$ws.Cells.Item(9, 4).Value = 0.6033474095103875

$ws.Cells.Item(10, 2).Value = 8782.122612132569
$ws.Cells.Item(10, 3).Value = 0.7254274587959104
$ws.Cells.Item(10, 4).Value = 0.5171030014022349

$ws.Cells.Item(11, 2).Value = 3748.111279331875
$ws.Cells.Item(11, 3).Value = 0.1059013586650847
$ws.Cells.Item(11, 4).Value = 0.2966265781769539

$ws.Cells.Item(12, 2).Value = 40684.3195025943
$ws.Cells.Item(12, 3).Value = 1011.027947386652
$ws.Cells.Item(12, 4).Value = 0.7680722285861548

$ws.Cells.Item(13, 2).Value = 59831.05914957315
$ws.Cells.Item(13, 3).Value = 133335.6317832718
$ws.Cells.Item(13, 4).Value = 0.7591197454263598

$ws.Cells.Item(14, 2).Value = 101861.6058959529
$ws.Cells.Item(14, 3).Value = 12422429266.50984
$ws.Cells.Item(14, 4).Value = 0.7426269218723918

$ws.Cells.Item(15, 2).Value = 118831.566280088
$ws.Cells.Item(15, 3).Value = 251607322218.1647
$ws.Cells.Item(15, 4).Value = 0.737861864990981

$ws.Cells.Item(16, 1).Value = "DM0 Subgroup 2"
$ws.Cells.Item(16, 2).Value = 90513.93200219698
$ws.Cells.Item(16, 3).Value = 34763481.07742978
$ws.Cells.Item(16, 4).Value = 0.8025489485682615

$ws.Cells.Item(17, 1).Value = "DM1 Subgroup 2"
$ws.Cells.Item(17, 2).Value = 101105.4581692971
$ws.Cells.Item(17, 3).Value = 235743195.4797715
$ws.Cells.Item(17, 4).Value = 0.7969355593891125

$ws.Cells.Item(18, 1).Value = "DM2 Subgroup 2"
$ws.Cells.Item(18, 2).Value = 105500.7629378873
$ws.Cells.Item(18, 3).Value = 153363870.480932
$ws.Cells.Item(18, 4).Value = 0.7946502055563922

$ws.Cells.Item(19, 2).Value = 36838.94176737932
$ws.Cells.Item(19, 3).Value = 736.3341818296014
$ws.Cells.Item(19, 4).Value = 0.7299976022606572

$ws.Cells.Item(20, 2).Value = 25342.88307628239
$ws.Cells.Item(20, 3).Value = 47.97312139119081
$ws.Cells.Item(20, 4).Value = 0.6958521485694259

$ws.Cells.Item(21, 2).Value = 13846.82438518547
$ws.Cells.Item(21, 3).Value = 2.482368678697976
$ws.Cells.Item(21, 4).Value = 0.6085472195774219

$ws.Cells.Item(22, 2).Value = -12505.63654077432
$ws.Cells.Item(22, 3).Value = -0.009898504804345587
$ws.Cells.Item(22, 4).Value = 0.9836595281802547

$ws.Cells.Item(23, 2).Value = -1836.207952824492
$ws.Cells.Item(23, 3).Value = -0.02513565892787702
$ws.Cells.Item(23, 4).Value = 0.4239329919827076

$ws.Cells.Item(24, 2).Value = 53508.52173692486
$ws.Cells.Item(24, 3).Value = 327525.1083798885
$ws.Cells.Item(24, 4).Value = 0.8161620435453979

$ws.Cells.Item(25, 2).Value = 75525.16978479302
$ws.Cells.Item(25, 3).Value = 18122315.57321637
$ws.Cells.Item(25, 4).Value = 0.8298872666862956

$ws.Cells.Item(26, 1).Value = "DM0 Subgroup 3"
$ws.Cells.Item(26, 2).Value = -11064.9001588159
$ws.Cells.Item(26, 3).Value = -0.01032559696440099
$ws.Cells.Item(26, 4).Value = 0.9657462999810005

$ws.Cells.Item(27, 1).Value = "DM1 Subgroup 3"
$ws.Cells.Item(27, 2).Value = -8676.21211520137
$ws.Cells.Item(27, 3).Value = -0.010403347419393
$ws.Cells.Item(27, 4).Value = 0.9264840029489525

$ws.Cells.Item(28, 1).Value = "DM2 Subgroup 3"
$ws.Cells.Item(28, 2).Value = -6659.121455476631
$ws.Cells.Item(28, 3).Value = -0.003172198020610486
$ws.Cells.Item(28, 4).Value = 0.8799100174518362

$ws.Cells.Item(29, 2).Value = -5872.929609460573
$ws.Cells.Item(29, 3).Value = -0.02942752314748293
$ws.Cells.Item(29, 4).Value = 0.9364025311519699

$ws.Cells.Item(30, 2).Value = -7891.290437778617
$ws.Cells.Item(30, 3).Value = -0.02392281016829638
$ws.Cells.Item(30, 4).Value = 0.9741242506812711

$ws.Cells.Item(31, 2).Value = -9909.651266096656
$ws.Cells.Item(31, 3).Value = -0.01817555953426603
$ws.Cells.Item(31, 4).Value = 0.9873993064886006
